$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'27.462.25"
$c.Style = "Normal"
$c = $ws.Range("E2")
$c.Value = "'  -2.29%  "
$c.Style = "Normal"
$c = $ws.Range("D3")
$c.Value = "'1.837.29"
$c.Style = "Normal"
$c = $ws.Range("E3")
$c.Value = "'  -2.86%  "
$c.Style = "Normal"
$c = $ws.Range("E4")
$c.Value = "'  -1.21%  "
$c.Style = "Normal"
$c = $ws.Range("D5")
$c.Value = "'331.82"
$c.Style = "Normal"
$c = $ws.Range("E5")
$c.Value = "'  -1.55%  "
$c.Style = "Normal"
$c = $ws.Range("E6")
$c.Value = "'  -1.09%  "
$c.Style = "Normal"
$c = $ws.Range("D7")
$c.Value = "'0.4609"
$c.Style = "Normal"
$c = $ws.Range("E7")
$c.Value = "'  -3.13%  "
$c.Style = "Normal"
$c = $ws.Range("E8")
$c.Value = "'  -3.73%  "
$c.Style = "Normal"
$c = $ws.Range("D9")
$c.Value = "'46.21"
$c.Style = "Normal"
$c = $ws.Range("E9")
$c.Value = "'  -2.19%  "
$c.Style = "Normal"
$c = $ws.Range("D10")
$c.Value = "'0.07864"
$c.Style = "Normal"
$c = $ws.Range("E10")
$c.Value = "'  -2.39%  "
$c.Style = "Normal"
$c = $ws.Range("D11")
$c.Value = "'0.9745"
$c.Style = "Normal"
$c = $ws.Range("E11")
$c.Value = "'  -4.77%  "
$c.Style = "Normal"
$c = $ws.Range("D12")
$c.Value = "'21.12"
$c.Style = "Normal"
$c = $ws.Range("E12")
$c.Value = "'  -4.08%  "
$c.Style = "Normal"
$c = $ws.Range("D13")
$c.Value = "'1.838.09"
$c.Style = "Normal"
$c = $ws.Range("E13")
$c.Value = "'  -1.84%  "
$c.Style = "Normal"
$c = $ws.Range("D14")
$c.Value = "'5.894"
$c.Style = "Normal"
$c = $ws.Range("E14")
$c.Value = "'  -2.53%  "
$c.Style = "Normal"
$c = $ws.Range("D15")
$c.Value = "'7.021"
$c.Style = "Normal"
$c = $ws.Range("E15")
$c.Value = "'  -3.09%  "
$c.Style = "Normal"
$c = $ws.Range("E16")
$c.Value = "'  -1.20%  "
$c.Style = "Normal"
$c = $ws.Range("D17")
$c.Value = "'87.86"
$c.Style = "Normal"
$c = $ws.Range("D18")
$c.Value = "'0.06639"
$c.Style = "Normal"
$c = $ws.Range("E18")
$c.Value = "'  -1.99%  "
$c.Style = "Normal"
$c = $ws.Range("D19")
$c.Value = "'0.00001030"
$c.Style = "Normal"
$c = $ws.Range("E19")
$c.Value = "'  -2.38%  "
$c.Style = "Normal"
$c = $ws.Range("E20")
$c.Value = "'  -1.17%  "
$c.Style = "Normal"
$c = $ws.Range("E21")
$c.Value = "'  -1.11%  "
$c.Style = "Normal"
$c = $ws.Range("D22")
$c.Value = "'27.461.74"
$c.Style = "Normal"
$c = $ws.Range("E22")
$c.Value = "'  -2.22%  "
$c.Style = "Normal"
$c = $ws.Range("D23")
$c.Value = "'5.336"
$c.Style = "Normal"
$c = $ws.Range("E23")
$c.Value = "'  -3.74%  "
$c.Style = "Normal"
$c = $ws.Range("D24")
$c.Value = "'10.85"
$c.Style = "Normal"
$c = $ws.Range("E24")
$c.Value = "'  -1.90%  "
$c.Style = "Normal"
$c = $ws.Range("D25")
$c.Value = "'2.296"
$c.Style = "Normal"
$c = $ws.Range("E25")
$c.Value = "'  -2.29%  "
$c.Style = "Normal"
$c = $ws.Range("D26")
$c.Value = "'157.13"
$c.Style = "Normal"
$c = $ws.Range("E26")
$c.Value = "'  -2.38%  "
$c.Style = "Normal"
$c = $ws.Range("D27")
$c.Value = "'19.35"
$c.Style = "Normal"
$c = $ws.Range("E27")
$c.Value = "'  -3.59%  "
$c.Style = "Normal"
$c = $ws.Range("D28")
$c.Value = "'2.066"
$c.Style = "Normal"
$c = $ws.Range("E28")
$c.Value = "'  -2.37%  "
$c.Style = "Normal"
$c = $ws.Range("D29")
$c.Value = "'5.326"
$c.Style = "Normal"
$c = $ws.Range("E29")
$c.Value = "'  -4.15%  "
$c.Style = "Normal"
$c = $ws.Range("D30")
$c.Value = "'118.62"
$c.Style = "Normal"
$c = $ws.Range("E30")
$c.Value = "'  -2.97%  "
$c.Style = "Normal"
$c = $ws.Range("E31")
$c.Value = "'  -2.90%  "
$c.Style = "Normal"
$c = $ws.Range("D32")
$c.Value = "'0.09292"
$c.Style = "Normal"
$c = $ws.Range("E32")
$c.Value = "'  -3.34%  "
$c.Style = "Normal"
$c = $ws.Range("E33")
$c.Value = "'  -1.98%  "
$c.Style = "Normal"
$c = $ws.Range("D34")
$c.Value = "'5.225"
$c.Style = "Normal"
$c = $ws.Range("E34")
$c.Value = "'  -2.82%  "
$c.Style = "Normal"
$c = $ws.Range("D35")
$c.Value = "'1.320"
$c.Style = "Normal"
$c = $ws.Range("E35")
$c.Value = "'  -3.80%  "
$c.Style = "Normal"
$c = $ws.Range("D36")
$c.Value = "'0.05934"
$c.Style = "Normal"
$c = $ws.Range("E36")
$c.Value = "'  -2.68%  "
$c.Style = "Normal"
$c = $ws.Range("D37")
$c.Value = "'0.02181"
$c.Style = "Normal"
$c = $ws.Range("E37")
$c.Value = "'  -3.49%  "
$c.Style = "Normal"
$c = $ws.Range("E38")
$c.Value = "'  -1.89%  "
$c.Style = "Normal"
$c = $ws.Range("D39")
$c.Value = "'1.157"
$c.Style = "Normal"
$c = $ws.Range("E39")
$c.Value = "'  -4.01%  "
$c.Style = "Normal"
$c = $ws.Range("D40")
$c.Value = "'0.5807"
$c.Style = "Normal"
$c = $ws.Range("E40")
$c.Value = "'  -3.05%  "
$c.Style = "Normal"
$c = $ws.Range("D41")
$c.Value = "'0.1838"
$c.Style = "Normal"
$c = $ws.Range("E41")
$c.Value = "'  -3.28%  "
$c.Style = "Normal"
$c = $ws.Range("D42")
$c.Value = "'10.05"
$c.Style = "Normal"
$c = $ws.Range("E42")
$c.Value = "'  -3.28%  "
$c.Style = "Normal"
$c = $ws.Range("D43")
$c.Value = "'1.238"
$c.Style = "Normal"
$c = $ws.Range("E43")
$c.Value = "'  -2.25%  "
$c.Style = "Normal"
$c = $ws.Range("E44")
$c.Value = "'  -1.23%  "
$c.Style = "Normal"
$c = $ws.Range("D45")
$c.Value = "'0.5484"
$c.Style = "Normal"
$c = $ws.Range("E45")
$c.Value = "'  -3.51%  "
$c.Style = "Normal"
$c = $ws.Range("D46")
$c.Value = "'1.865"
$c.Style = "Normal"
$c = $ws.Range("E46")
$c.Value = "'  -3.82%  "
$c.Style = "Normal"
$c = $ws.Range("D47")
$c.Value = "'0.06652"
$c.Style = "Normal"
$c = $ws.Range("E47")
$c.Value = "'  -2.71%  "
$c.Style = "Normal"
$c = $ws.Range("D48")
$c.Value = "'109.77"
$c.Style = "Normal"
$c = $ws.Range("E48")
$c.Value = "'  -2.49%  "
$c.Style = "Normal"
$c = $ws.Range("E49")
$c.Value = "'  -3.25%  "
$c.Style = "Normal"
$c = $ws.Range("B50")
$c.Value = "'PaxDollar"
$c.Style = "Normal"
$c = $ws.Range("C50")
$c.Value = "'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$c.Style = "Normal"
$c = $ws.Range("D50")
$c.Value = "'1.001"
$c.Style = "Normal"
$c = $ws.Range("E50")
$c.Value = "'  -1.26%  "
$c.Style = "Normal"
$c = $ws.Range("B51")
$c.Value = "'BabyDogeCoin"
$c.Style = "Normal"
$c = $ws.Range("C51")
$c.Value = "'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$c.Style = "Normal"
$c = $ws.Range("D51")
$c.Value = "'0.00000000288"
$c.Style = "Normal"
$c = $ws.Range("E51")
$c.Value = "'  -1.00%  "
$c.Style = "Normal"
